$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-31 hold three regions' (Brasil, Nordeste, Sergipe) yearly values.
# The series was re-pointed three years earlier (e.g. 2014 -> 2011) and the
# "Valor" figures were refreshed to match the newly-aligned years.
$updates = @(
    @{ Row = 2; B = "01/01/2011"; C = 2199 },
    @{ Row = 3; B = "01/01/2012"; C = 2258.15 },
    @{ Row = 4; B = "01/01/2013"; C = 2316.99 },
    @{ Row = 5; B = "01/01/2014"; C = 2356.25 },
    @{ Row = 6; B = "01/01/2015"; C = 2291.66 },
    @{ Row = 7; B = "01/01/2016"; C = 2261.01 },
    @{ Row = 8; B = "01/01/2017"; C = 2272.25 },
    @{ Row = 9; B = "01/01/2018"; C = 2298.54 },
    @{ Row = 10; B = "01/01/2019"; C = 2323.85 },
    @{ Row = 11; B = "01/01/2020"; C = 2273.25 },
    @{ Row = 12; B = "01/01/2011"; C = 1328.31 },
    @{ Row = 13; B = "01/01/2012"; C = 1388.83 },
    @{ Row = 14; B = "01/01/2013"; C = 1456.43 },
    @{ Row = 15; B = "01/01/2014"; C = 1468.71 },
    @{ Row = 16; B = "01/01/2015"; C = 1432.96 },
    @{ Row = 17; B = "01/01/2016"; C = 1431.13 },
    @{ Row = 18; B = "01/01/2017"; C = 1428.68 },
    @{ Row = 19; B = "01/01/2018"; C = 1438.08 },
    @{ Row = 20; B = "01/01/2019"; C = 1475.58 },
    @{ Row = 21; B = "01/01/2020"; C = 1434.5 },
    @{ Row = 22; B = "01/01/2011"; C = 1640.7 },
    @{ Row = 23; B = "01/01/2012"; C = 1695.4 },
    @{ Row = 24; B = "01/01/2013"; C = 1774.82 },
    @{ Row = 25; B = "01/01/2014"; C = 1785.98 },
    @{ Row = 26; B = "01/01/2015"; C = 1756.28 },
    @{ Row = 27; B = "01/01/2016"; C = 1715.48 },
    @{ Row = 28; B = "01/01/2017"; C = 1652.14 },
    @{ Row = 29; B = "01/01/2018"; C = 1674.45 },
    @{ Row = 30; B = "01/01/2019"; C = 1655.26 },
    @{ Row = 31; B = "01/01/2020"; C = 1611.88 }
)

foreach ($u in $updates) {
    $bCell = $ws.Range("B" + $u.Row)
    # Format as text first so Excel keeps the literal "dd/mm/yyyy" string
    # instead of auto-converting it to a date serial number.
    $bCell.NumberFormat = "@"
    $bCell.Value = $u.B
    $bCell.Style = "Normal"

    $ws.Range("C" + $u.Row).Value = $u.C
}
